$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2742.3276
$ws.Range("I137").Value = 2929.64
$ws.Range("J137").Value = 2600.4243
$ws.Range("K137").Value = 8788.92
$ws.Range("L137").Value = 7801.2729
$ws.Range("M137").Value = -6238.92
$ws.Range("N137").Value = -12901.2729
$ws.Range("H138").Value = 8866.65
$ws.Range("J138").Value = 9433.625
$ws.Range("L138").Value = 28300.875
$ws.Range("N138").Value = -38580.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 749
$ws.Range("J5").Value = 749
$ws.Range("L5").Value = 749
$ws.Range("N5").Value = -973
$ws.Range("H44").Value = 89999
$ws.Range("J44").Value = 89999
$ws.Range("L44").Value = 89999
$ws.Range("N44").Value = -90975
$ws.Range("H55").Value = 68633.336
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 68633.336
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 68633.336
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -69263.336
$ws.Range("H80").Value = 82191.28999999999
$ws.Range("J80").Value = 82191.28999999999
$ws.Range("L80").Value = 82191.28999999999
$ws.Range("N80").Value = -84187.28999999999
$ws.Range("H82").Value = 29931.334
$ws.Range("I82").Value = 30000
$ws.Range("J82").Value = 29897
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 29897
$ws.Range("M82").Value = -29639
$ws.Range("N82").Value = -30619
$ws.Range("H83").Value = 82191.28999999999
$ws.Range("J83").Value = 82191.28999999999
$ws.Range("L83").Value = 246573.87
$ws.Range("N83").Value = -256557.87
$ws.Range("H85").Value = 29931.334
$ws.Range("I85").Value = 30000
$ws.Range("J85").Value = 29897
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 29897
$ws.Range("M85").Value = -28752
$ws.Range("N85").Value = -32393
$ws.Range("H97").Value = 1522.5769
$ws.Range("I97").Value = 1522.5769
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1522.5769
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1026.5769
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 3632.2666
$ws.Range("I122").Value = 4262
$ws.Range("J122").Value = 3562.2964
$ws.Range("K122").Value = 12786
$ws.Range("L122").Value = 10686.8892
$ws.Range("M122").Value = -10336
$ws.Range("N122").Value = -15586.8892
$ws.Range("H123").Value = 79247.5
$ws.Range("J123").Value = 79247.5
$ws.Range("L123").Value = 79247.5
$ws.Range("N123").Value = -89047.5
$ws.Range("H132").Value = 5008.1577
$ws.Range("I132").Value = 3263
$ws.Range("K132").Value = 9789
$ws.Range("M132").Value = -7259

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 749
$ws.Range("J4").Value = 749
$ws.Range("L4").Value = 749
$ws.Range("N4").Value = -979
$ws.Range("H20").Value = 2028.6471
$ws.Range("I20").Value = 1747.3636
$ws.Range("J20").Value = 2544.3333
$ws.Range("K20").Value = 1747.3636
$ws.Range("L20").Value = 2544.3333
$ws.Range("M20").Value = -1500.3636
$ws.Range("N20").Value = -3038.3333
$ws.Range("H22").Value = 263.5
$ws.Range("I22").Value = 239.5
$ws.Range("K22").Value = 239.5
$ws.Range("M22").Value = -66.5
$ws.Range("H35").Value = 72671.28999999999
$ws.Range("J35").Value = 72671.28999999999
$ws.Range("L35").Value = 72671.28999999999
$ws.Range("N35").Value = -73291.28999999999
$ws.Range("H80").Value = 1010.6667
$ws.Range("J80").Value = 762.1429000000001
$ws.Range("L80").Value = 762.1429000000001
$ws.Range("N80").Value = -2758.1429
$ws.Range("H82").Value = 31006.268
$ws.Range("I82").Value = 10629.4
$ws.Range("J82").Value = 71760
$ws.Range("K82").Value = 10629.4
$ws.Range("L82").Value = 71760
$ws.Range("M82").Value = -10246.4
$ws.Range("N82").Value = -72526
$ws.Range("H83").Value = 1010.6667
$ws.Range("J83").Value = 762.1429000000001
$ws.Range("L83").Value = 3810.7145
$ws.Range("N83").Value = -13794.7145
$ws.Range("H85").Value = 31006.268
$ws.Range("I85").Value = 10629.4
$ws.Range("J85").Value = 71760
$ws.Range("K85").Value = 10629.4
$ws.Range("L85").Value = 71760
$ws.Range("M85").Value = -9303.4
$ws.Range("N85").Value = -74412
$ws.Range("H99").Value = 3287.3
$ws.Range("I99").Value = 2340
$ws.Range("K99").Value = 2340
$ws.Range("M99").Value = -842
$ws.Range("H105").Value = 2653.45
$ws.Range("I105").Value = 1975.1538
$ws.Range("J105").Value = 3913.1428
$ws.Range("K105").Value = 1975.1538
$ws.Range("L105").Value = 3913.1428
$ws.Range("M105").Value = -228.1538
$ws.Range("N105").Value = -7407.1428
$ws.Range("H116").Value = 72749.5
$ws.Range("J116").Value = 72749.5
$ws.Range("L116").Value = 72749.5
$ws.Range("N116").Value = -81927.5
$ws.Range("H117").ClearContents()
$ws.Range("I117").ClearContents()
$ws.Range("J117").ClearContents()
$ws.Range("K117").ClearContents()
$ws.Range("L117").ClearContents()
$ws.Range("H118").ClearContents()
$ws.Range("I118").ClearContents()
$ws.Range("J118").ClearContents()
$ws.Range("K118").ClearContents()
$ws.Range("L118").ClearContents()
$ws.Range("N118").ClearContents()
$ws.Range("H119").ClearContents()
$ws.Range("I119").ClearContents()
$ws.Range("J119").ClearContents()
$ws.Range("K119").ClearContents()
$ws.Range("L119").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("H120").ClearContents()
$ws.Range("I120").ClearContents()
$ws.Range("J120").ClearContents()
$ws.Range("K120").ClearContents()
$ws.Range("L120").ClearContents()
$ws.Range("H122").ClearContents()
$ws.Range("I122").ClearContents()
$ws.Range("J122").ClearContents()
$ws.Range("K122").ClearContents()
$ws.Range("L122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H123").ClearContents()
$ws.Range("I123").ClearContents()
$ws.Range("J123").ClearContents()
$ws.Range("K123").ClearContents()
$ws.Range("L123").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("M128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2337.125
$ws.Range("I16").Value = 1790
$ws.Range("J16").Value = 3978.5
$ws.Range("K16").Value = 1790
$ws.Range("L16").Value = 3978.5
$ws.Range("M16").Value = -1503
$ws.Range("N16").Value = -4552.5
$ws.Range("H22").Value = 396.16666
$ws.Range("I22").Value = 286.25
$ws.Range("K22").Value = 286.25
$ws.Range("M22").Value = 63.75
$ws.Range("H31").Value = 3584.7942
$ws.Range("I31").Value = 1014.06665
$ws.Range("K31").Value = 1014.06665
$ws.Range("M31").Value = -719.06665
$ws.Range("H34").Value = 3584.7942
$ws.Range("I34").Value = 1014.06665
$ws.Range("K34").Value = 1014.06665
$ws.Range("M34").Value = -812.06665
$ws.Range("H50").Value = 63895.75
$ws.Range("J50").Value = 63895.75
$ws.Range("L50").Value = 63895.75
$ws.Range("N50").Value = -65145.75
$ws.Range("H58").Value = 5551.7646
$ws.Range("I58").Value = 3412.1428
$ws.Range("J58").Value = 7049.5
$ws.Range("K58").Value = 3412.1428
$ws.Range("L58").Value = 7049.5
$ws.Range("M58").Value = -3209.1428
$ws.Range("N58").Value = -7455.5
$ws.Range("H60").Value = 77167.336
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 77167.336
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 77167.336
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -78189.336
$ws.Range("H62").Value = 2730.25
$ws.Range("I62").Value = 2730.25
$ws.Range("K62").Value = 2730.25
$ws.Range("M62").Value = -2106.25
$ws.Range("H65").Value = 2730.25
$ws.Range("I65").Value = 2730.25
$ws.Range("K65").Value = 13651.25
$ws.Range("M65").Value = -10531.25
$ws.Range("H68").Value = 76259
$ws.Range("J68").Value = 87823.75
$ws.Range("L68").Value = 87823.75
$ws.Range("N68").Value = -89321.75
$ws.Range("H69").Value = 44996.668
$ws.Range("I69").Value = 55000
$ws.Range("J69").Value = 24990
$ws.Range("K69").Value = 55000
$ws.Range("L69").Value = 24990
$ws.Range("M69").Value = -54251
$ws.Range("N69").Value = -26488
$ws.Range("H71").Value = 76259
$ws.Range("J71").Value = 87823.75
$ws.Range("L71").Value = 263471.25
$ws.Range("N71").Value = -270959.25
$ws.Range("H72").Value = 44996.668
$ws.Range("I72").Value = 55000
$ws.Range("J72").Value = 24990
$ws.Range("K72").Value = 165000
$ws.Range("L72").Value = 74970
$ws.Range("M72").Value = -161256
$ws.Range("N72").Value = -82458
$ws.Range("H99").Value = 6555.3335
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 7333
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 7333
$ws.Range("M99").Value = -3502
$ws.Range("N99").Value = -10329
$ws.Range("H113").Value = 2337.125
$ws.Range("I113").Value = 1790
$ws.Range("J113").Value = 3978.5
$ws.Range("K113").Value = 1790
$ws.Range("L113").Value = 3978.5
$ws.Range("M113").Value = 380
$ws.Range("N113").Value = -8318.5
$ws.Range("H126").Value = 6555.3335
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 7333
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 21999
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -26939
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 58992.11
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 58992.11
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 58992.11
$ws.Range("N130").Value = -69032.11
$ws.Range("H131").Value = 10296
$ws.Range("I131").Value = 10296
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 10296
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -5256
$ws.Range("H132").Value = 3929.7083
$ws.Range("I132").Value = 3999.8948
$ws.Range("J132").Value = 3663
$ws.Range("K132").Value = 11999.6844
$ws.Range("L132").Value = 10989
$ws.Range("M132").Value = -9469.6844
$ws.Range("N132").Value = -16049
$ws.Range("H133").Value = 58571.145
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 58571.145
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 58571.145
$ws.Range("N133").Value = -63631.145
$ws.Range("H134").Value = 273604.28
$ws.Range("I134").Value = 2561.5186
$ws.Range("J134").Value = 1005419.7
$ws.Range("K134").Value = 7684.5558
$ws.Range("L134").Value = 3016259.1
$ws.Range("M134").Value = -5149.5558
$ws.Range("N134").Value = -3021329.1
$ws.Range("H135").Value = 74800.234
$ws.Range("I135").Value = 70700
$ws.Range("J135").Value = 75005.25
$ws.Range("K135").Value = 70700
$ws.Range("L135").Value = 75005.25
$ws.Range("M135").Value = -65630
$ws.Range("N135").Value = -85145.25
$ws.Range("H136").Value = 5551.7646
$ws.Range("I136").Value = 3412.1428
$ws.Range("J136").Value = 7049.5
$ws.Range("K136").Value = 10236.4284
$ws.Range("L136").Value = 21148.5
$ws.Range("M136").Value = -7686.428400000001
$ws.Range("N136").Value = -26248.5
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 220280
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 220280
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 220280
$ws.Range("N138").Value = -230560
$ws.Range("H139").Value = 92548
$ws.Range("I139").Value = 95700
$ws.Range("J139").Value = 90446.664
$ws.Range("K139").Value = 95700
$ws.Range("L139").Value = 90446.664
$ws.Range("M139").Value = -90560
$ws.Range("N139").Value = -100726.664
$ws.Range("H140").Value = 95999.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 95999.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 95999.8
$ws.Range("N140").Value = -106359.8
$ws.Range("H141").Value = 109031.1
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 109031.1
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 109031.1
$ws.Range("N141").Value = -119391.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 9083.75
$ws.Range("I14").Value = 9083.75
$ws.Range("K14").Value = 27251.25
$ws.Range("M14").Value = -27078.25
$ws.Range("H64").Value = 4003
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 4337.6665
$ws.Range("K64").Value = 8997
$ws.Range("L64").Value = 13012.9995
$ws.Range("M64").Value = -8727
$ws.Range("N64").Value = -13552.9995
$ws.Range("H67").Value = 4003
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 4337.6665
$ws.Range("K67").Value = 8997
$ws.Range("L67").Value = 13012.9995
$ws.Range("M67").Value = -8061
$ws.Range("N67").Value = -14884.9995
$ws.Range("H68").Value = 2303.111
$ws.Range("J68").Value = 2509.7917
$ws.Range("L68").Value = 7529.375100000001
$ws.Range("N68").Value = -9151.375100000001
$ws.Range("H70").Value = 251497.5
$ws.Range("I70").Value = 251497.5
$ws.Range("K70").Value = 754492.5
$ws.Range("M70").Value = -754177.5
$ws.Range("H71").Value = 2303.111
$ws.Range("J71").Value = 2509.7917
$ws.Range("L71").Value = 22588.1253
$ws.Range("N71").Value = -30700.1253
$ws.Range("H73").Value = 251497.5
$ws.Range("I73").Value = 251497.5
$ws.Range("K73").Value = 754492.5
$ws.Range("M73").Value = -753400.5
$ws.Range("H107").Value = 3166.625
$ws.Range("J107").Value = 4935.1665
$ws.Range("L107").Value = 14805.4995
$ws.Range("N107").Value = -18645.4995
$ws.Range("H117").Value = 420.73334
$ws.Range("I117").Value = 453.375
$ws.Range("J117").Value = 383.42856
$ws.Range("K117").Value = 1360.125
$ws.Range("L117").Value = 1150.28568
$ws.Range("M117").Value = 2081.875
$ws.Range("N117").Value = -8034.28568
$ws.Range("H140").Value = 2443.7646
$ws.Range("I140").Value = 2443.7646
$ws.Range("K140").Value = 7331.293799999999
$ws.Range("M140").Value = -2151.293799999999
$ws.Range("H141").Value = 1246.2858
$ws.Range("I141").Value = 1246.2858
$ws.Range("K141").Value = 3738.8574
$ws.Range("M141").Value = 1441.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3587
$ws.Range("I99").Value = 3587
$ws.Range("K99").Value = 3587
$ws.Range("M99").Value = -1341
$ws.Range("H113").Value = 602579.75
$ws.Range("I113").Value = 1127317.4
$ws.Range("J113").Value = 12250
$ws.Range("K113").Value = 1127317.4
$ws.Range("L113").Value = 12250
$ws.Range("M113").Value = -1125147.4
$ws.Range("N113").Value = -16590
$ws.Range("H121").Value = 57072.5
$ws.Range("J121").Value = 57072.5
$ws.Range("L121").Value = 57072.5
$ws.Range("N121").Value = -60566.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 30
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 30
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 380
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 3339574
$ws.Range("I40").Value = 4291188.5
$ws.Range("J40").Value = 8923
$ws.Range("K40").Value = 4291188.5
$ws.Range("L40").Value = 8923
$ws.Range("M40").Value = -4291052.5
$ws.Range("N40").Value = -9195
$ws.Range("H46").Value = 3534.5938
$ws.Range("I46").Value = 2950.5
$ws.Range("J46").Value = 4508.0835
$ws.Range("K46").Value = 2950.5
$ws.Range("L46").Value = 4508.0835
$ws.Range("M46").Value = -2762.5
$ws.Range("N46").Value = -4884.0835
$ws.Range("H55").Value = 579
$ws.Range("I55").Value = 329.9375
$ws.Range("K55").Value = 329.9375
$ws.Range("M55").Value = -156.9375
$ws.Range("H61").Value = 7811.1113
$ws.Range("I61").Value = 6836.364
$ws.Range("K61").Value = 6836.364
$ws.Range("M61").Value = -6634.364
$ws.Range("H68").Value = 6186.6665
$ws.Range("J68").Value = 5896.5
$ws.Range("L68").Value = 5896.5
$ws.Range("N68").Value = -7394.5
$ws.Range("H71").Value = 6186.6665
$ws.Range("J71").Value = 5896.5
$ws.Range("L71").Value = 29482.5
$ws.Range("N71").Value = -36970.5
$ws.Range("H113").Value = 7811.1113
$ws.Range("I113").Value = 6836.364
$ws.Range("K113").Value = 6836.364
$ws.Range("M113").Value = -4666.364
$ws.Range("H136").Value = 4118.0215
$ws.Range("I136").Value = 3216.5833
$ws.Range("J136").Value = 7068.1816
$ws.Range("K136").Value = 9649.749899999999
$ws.Range("L136").Value = 21204.5448
$ws.Range("M136").Value = -7099.749899999999
$ws.Range("N136").Value = -26304.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10365.111
$ws.Range("J45").Value = 10660.75
$ws.Range("L45").Value = 10660.75
$ws.Range("N45").Value = -11642.75
$ws.Range("H107").Value = 24981.453
$ws.Range("I107").Value = 41426
$ws.Range("J107").Value = 798.2941
$ws.Range("K107").Value = 124278
$ws.Range("L107").Value = 2394.8823
$ws.Range("M107").Value = -122358
$ws.Range("N107").Value = -6234.882299999999
$ws.Range("H122").Value = 3692.158
$ws.Range("I122").Value = 2471.125
$ws.Range("K122").Value = 7413.375
$ws.Range("M122").Value = -4963.375
$ws.Range("H126").Value = 2663.7827
$ws.Range("I126").Value = 2413.35
$ws.Range("K126").Value = 7240.049999999999
$ws.Range("M126").Value = -4770.049999999999
$ws.Range("H136").Value = 449471.66
$ws.Range("J136").Value = 161329.16
$ws.Range("L136").Value = 483987.48
$ws.Range("N136").Value = -489087.48
